# Group1 Test Case Plan - add Sam's tests (Functional Tests section) and
# introduce a "Tester" column, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: insert a new "Tester" column between "Test Case
#    Description" (E) and "Status" (old F). The former spacer column
#    (old G, hidden/width 0) becomes the new "Status" column, and a new
#    narrow hidden column is used for "Tester".
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Tester"
$ws.Range("G1").Value = "Status"

# ---------------------------------------------------------------------
# 2. Move the "Not started" status values that used to live in column F
#    (rows 2-8) over to column G, and blank out column F for those rows.
# ---------------------------------------------------------------------
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Not started"
$ws.Range("G2").WrapText = $true

$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "Not started"
$ws.Range("G3").WrapText = $true
$ws.Range("E3").WrapText = $true

$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Not started"
$ws.Range("G4").WrapText = $true

$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "Not started"
$ws.Range("G5").WrapText = $true

$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = "Not started"
$ws.Range("G6").WrapText = $true

$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = "Not started"
$ws.Range("G7").WrapText = $true

$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

# ---------------------------------------------------------------------
# 3. Column widths / visibility: E grows a touch and wraps, F becomes
#    the new narrow hidden "Tester" column, G takes over the old
#    "Status" column width and becomes visible.
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 33.166666666666664
$ws.Columns("F").ColumnWidth = 8
$ws.Columns("F").Hidden = $true
$ws.Columns("G").Hidden = $false
$ws.Columns("G").ColumnWidth = 13.333333333333334

# ---------------------------------------------------------------------
# 4. Rows 2-7 (the old scenario rows) are now hidden / collapsed.
# ---------------------------------------------------------------------
$ws.Rows("2:7").Hidden = $true

# ---------------------------------------------------------------------
# 5. New "Sprint 2" / "Functional Tests" section header (row 9).
# ---------------------------------------------------------------------
$ws.Range("B9").Value = "Sprint 2"
$ws.Range("D9").Value = "Functional Tests"

# ---------------------------------------------------------------------
# 6. Sam's new functional test rows (11-15), leaving row 10 blank.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "testNewRacer"
$ws.Range("E11").Value = "Test creation of a Racer and associated getters"
$ws.Range("E11").WrapText = $true
$ws.Range("F11").Value = "Sam"
$ws.Range("G11").Value = "Ran - Passed"
$ws.Rows("11").RowHeight = 28.8

$ws.Range("D12").Value = "testRacerStartFinish"
$ws.Range("E12").Value = "Test starting and stoping a Racer and associated getters "
$ws.Range("E12").WrapText = $true
$ws.Range("F12").Value = "Sam"
$ws.Range("G12").Value = "Ran - Passed"
$ws.Rows("12").RowHeight = 28.8

$ws.Range("D13").Value = "testNewStartSensor"
$ws.Range("E13").Value = "Test adding a Sensor to a start trigger"
$ws.Range("E13").WrapText = $true
$ws.Range("F13").Value = "Sam"
$ws.Range("G13").Value = "Ran - Passed"
$ws.Range("H13").Value = "merge with stopsensor?"

$ws.Range("D14").Value = "testNewStopSensor"
$ws.Range("E14").Value = "Test adding a Sensor to a stop trigger"
$ws.Range("E14").WrapText = $true
$ws.Range("F14").Value = "Sam"
$ws.Range("G14").Value = "Ran - Passed"
$ws.Range("H14").Value = "merge with startsensor?"

$ws.Range("D15").Value = "TestToggleSensor"
$ws.Range("E15").Value = "Test ability to toggle a Sensor on and off"
$ws.Range("E15").WrapText = $true
$ws.Range("F15").Value = "Sam"
$ws.Range("G15").Value = "Ran - Passed"
$ws.Range("H15").Value = "merge with other sensors?"
$ws.Rows("15").RowHeight = 28.8

# ---------------------------------------------------------------------
# 7. Selection / view state.
# ---------------------------------------------------------------------
$ws.Range("A2:XFD7").Select()
